# Update automàtic: dades i banners [2026-02-20 10:47]
#
# Re-runs the meteo.cat XEMA scrape for station XJ: the 09:30-10:00 UTC
# period row is refreshed to the 10:00-10:30 UTC period (new readings +
# new extraction timestamps), and the extraction timestamp on the three
# "yesterday" rows is bumped to match the new run. The header-discovery
# sheet's sample URL is refreshed to match.

$wb = $excel.ActiveWorkbook

# Worksheets.Item by index avoids any encoding headaches with the
# accented sheet names ("Dades_Període" / "Estudi_Capçaleres").
$wsDades = $wb.Worksheets.Item(1)
$wsCap   = $wb.Worksheets.Item(2)

# Helper: write a literal string into a cell without Excel's COM layer
# re-typing numeric-looking text ("184", "1024.5", ...) as a Number.
# Forcing the number format to Text ("@") before the assignment makes
# the write land as a string; resetting the style back to "Normal"
# afterwards drops the format override again so no stray cell styling
# is left behind.
function Set-TextValue {
    param(
        $Range,
        [string]$Value
    )
    $Range.NumberFormat = "@"
    $Range.Value2 = $Value
    $Range.Style = "Normal"
}

# ---- Dades_Període, row 2 (the 09:30-10:00 -> 10:00-10:30 period) ----

Set-TextValue $wsDades.Range("E2")  "10:00 - 10:30"
Set-TextValue $wsDades.Range("H2")  "2026-02-20 10:47:34"
Set-TextValue $wsDades.Range("I2")  "10:00"
Set-TextValue $wsDades.Range("J2")  "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T10:00Z"

Set-TextValue $wsDades.Range("M2")  "184"
Set-TextValue $wsDades.Range("N2")  "46"
Set-TextValue $wsDades.Range("O2")  "1024.5"

Set-TextValue $wsDades.Range("Q2")  "10:00 - 10:30"
Set-TextValue $wsDades.Range("R2")  "414"
Set-TextValue $wsDades.Range("S2")  "13.6"
Set-TextValue $wsDades.Range("T2")  "13.0"
Set-TextValue $wsDades.Range("U2")  "14.4"
Set-TextValue $wsDades.Range("V2")  "2.2"
Set-TextValue $wsDades.Range("W2")  "8.6"

Set-TextValue $wsDades.Range("X2")  "10:00 - 10:30"
Set-TextValue $wsDades.Range("Y2")  "13.6"
Set-TextValue $wsDades.Range("Z2")  "14.4"
Set-TextValue $wsDades.Range("AA2") "13.0"
Set-TextValue $wsDades.Range("AB2") "46"
Set-TextValue $wsDades.Range("AD2") "2.2"
Set-TextValue $wsDades.Range("AE2") "184"
Set-TextValue $wsDades.Range("AF2") "8.6"
Set-TextValue $wsDades.Range("AG2") "1024.5"
Set-TextValue $wsDades.Range("AH2") "414"

Set-TextValue $wsDades.Range("AI2") "10:00 - 10:30"
Set-TextValue $wsDades.Range("AJ2") "13.6"
Set-TextValue $wsDades.Range("AK2") "14.4"
Set-TextValue $wsDades.Range("AL2") "13.0"
Set-TextValue $wsDades.Range("AM2") "46"
Set-TextValue $wsDades.Range("AO2") "2.2"
Set-TextValue $wsDades.Range("AP2") "184"
Set-TextValue $wsDades.Range("AQ2") "8.6"
Set-TextValue $wsDades.Range("AR2") "1024.5"
Set-TextValue $wsDades.Range("AS2") "414"

# ---- Dades_Període, rows 3-6: extraction timestamp bump only ----

Set-TextValue $wsDades.Range("H3") "2026-02-20 10:47:36"
Set-TextValue $wsDades.Range("H4") "2026-02-20 10:47:36"
Set-TextValue $wsDades.Range("H5") "2026-02-20 10:47:36"
Set-TextValue $wsDades.Range("H6") "2026-02-20 10:47:36"

# ---- Estudi_Capçaleres, row 2: refreshed sample source URL ----

Set-TextValue $wsCap.Range("F2") "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T10:00Z"
